$d = $word.ActiveDocument

# Find the paragraph containing "Docente(s) Responsável(eis)"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Docente(s) Responsável(eis)*") {
        $target = $p
        break
    }
}

# Insert a new paragraph right after the target paragraph
$insertRange = $target.Range
$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()

# The newly inserted paragraph is the next one after $target
$newPara = $target.Next()
$newPara.Range.Text = "5701460 - Antonio Iacono"
$newPara.Style = "ListBullet"
